# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 314
$ws1.Range("F4").Value = 238
$ws1.Range("F5").Value = 2788
$ws1.Range("F6").Value = 1953
$ws1.Range("F8").Value = 124
$ws1.Range("F9").Value = 1008
$ws1.Range("F10").Value = 189
$ws1.Range("F11").Value = 43
$ws1.Range("F12").Value = 30

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 314
$ws4.Range("F4").Value = 238
$ws4.Range("F5").Value = 2788
$ws4.Range("F6").Value = 1953
$ws4.Range("F9").Value = 124
$ws4.Range("F10").Value = 1008
$ws4.Range("F11").Value = 189
$ws4.Range("F12").Value = 43
$ws4.Range("F13").Value = 30
